$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1498
$ws.Range("I32").Value = 995
$ws.Range("K32").Value = 995
$ws.Range("M32").Value = -669

$ws.Range("H43").Value = 1184.8572
$ws.Range("J43").Value = 1423.5
$ws.Range("L43").Value = 1423.5
$ws.Range("N43").Value = -1561.5

$ws.Range("H70").Value = 63808.875
$ws.Range("I70").Value = 167717
$ws.Range("J70").Value = 1464
$ws.Range("K70").Value = 503151
$ws.Range("L70").Value = 4392
$ws.Range("M70").Value = -502881
$ws.Range("N70").Value = -4932

$ws.Range("H73").Value = 63808.875
$ws.Range("I73").Value = 167717
$ws.Range("J73").Value = 1464
$ws.Range("K73").Value = 503151
$ws.Range("L73").Value = 4392
$ws.Range("M73").Value = -502215
$ws.Range("N73").Value = -6264

$ws.Range("H98").Value = 3546
$ws.Range("I98").Value = 2516.1538
$ws.Range("K98").Value = 2516.1538
$ws.Range("M98").Value = -1018.1538

$ws.Range("H112").Value = 1843.8182
$ws.Range("J112").Value = 1948.2
$ws.Range("L112").Value = 5844.6
$ws.Range("N112").Value = -8060.6

$ws.Range("H122").Value = 3546
$ws.Range("I122").Value = 2516.1538
$ws.Range("K122").Value = 7548.4614
$ws.Range("M122").Value = -5098.4614

$ws.Range("H137").Value = 1481.1875
$ws.Range("I137").Value = 1315.68
$ws.Range("K137").Value = 3947.04
$ws.Range("M137").Value = -1397.04

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H27").Value = 10000
$ws.Range("J27").Value = 10000
$ws.Range("L27").Value = 10000
$ws.Range("N27").Value = -10368

$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()

$ws.Range("H52").Value = 26794.75
$ws.Range("J52").Value = 26794.75
$ws.Range("L52").Value = 26794.75
$ws.Range("N52").Value = -27430.75

$ws.Range("H74").Value = 1516.8
$ws.Range("I74").Value = 1309.7142
$ws.Range("K74").Value = 1309.7142
$ws.Range("M74").Value = -435.7141999999999

$ws.Range("H77").Value = 1516.8
$ws.Range("I77").Value = 1309.7142
$ws.Range("K77").Value = 6548.571
$ws.Range("M77").Value = -2180.571

$ws.Range("H92").Value = 135593.8
$ws.Range("J92").Value = 135593.8
$ws.Range("L92").Value = 135593.8
$ws.Range("N92").Value = -140585.8

$ws.Range("H101").Value = 42850.75
$ws.Range("J101").Value = 42850.75
$ws.Range("L101").Value = 42850.75
$ws.Range("N101").Value = -49340.75

$ws.Range("H123").Value = 29940.8
$ws.Range("J123").Value = 29940.8
$ws.Range("L123").Value = 29940.8
$ws.Range("N123").Value = -39740.8

$ws.Range("H132").Value = 1717.3334
$ws.Range("I132").Value = 1343.5714
$ws.Range("J132").Value = 2723.6155
$ws.Range("K132").Value = 4030.7142
$ws.Range("L132").Value = 8170.8465
$ws.Range("M132").Value = -1500.7142
$ws.Range("N132").Value = -13230.8465

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4169
$ws.Range("I105").Value = 3987.1428
$ws.Range("J105").Value = 4593.3335
$ws.Range("K105").Value = 3987.1428
$ws.Range("L105").Value = 4593.3335
$ws.Range("M105").Value = -2240.1428
$ws.Range("N105").Value = -8087.3335

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 712.3333
$ws.Range("I2").Value = 1338.625
$ws.Range("J2").Value = 326.92307
$ws.Range("K2").Value = 8031.75
$ws.Range("L2").Value = 1961.53842
$ws.Range("M2").Value = -7918.75
$ws.Range("N2").Value = -2187.53842

$ws.Range("H9").Value = 500200.5
$ws.Range("I9").Value = 400
$ws.Range("K9").Value = 1200
$ws.Range("M9").Value = -976

$ws.Range("H22").Value = 20834976
$ws.Range("J22").Value = 1877.1428
$ws.Range("L22").Value = 5631.428400000001
$ws.Range("N22").Value = -5969.428400000001

$ws.Range("H23").Value = 411.25
$ws.Range("J23").Value = 398.57144
$ws.Range("L23").Value = 1195.71432
$ws.Range("N23").Value = -1665.71432

$ws.Range("H27").Value = 20834976
$ws.Range("J27").Value = 1877.1428
$ws.Range("L27").Value = 5631.428400000001
$ws.Range("N27").Value = -5835.428400000001

$ws.Range("H58").Value = 4533.3335
$ws.Range("J58").Value = 4533.3335
$ws.Range("L58").Value = 13600.0005
$ws.Range("N58").Value = -13856.0005

$ws.Range("H107").Value = 655.63635
$ws.Range("I107").Value = 920
$ws.Range("K107").Value = 2760
$ws.Range("M107").Value = -840

$ws.Range("H131").Value = 14087195
$ws.Range("I131").Value = 10355
$ws.Range("J131").Value = 16394873
$ws.Range("K131").Value = 31065
$ws.Range("L131").Value = 49184619
$ws.Range("M131").Value = -26025
$ws.Range("N131").Value = -49194699

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 34400
$ws.Range("J63").Value = 34400
$ws.Range("L63").Value = 34400
$ws.Range("N63").Value = -35772

$ws.Range("H66").Value = 34400
$ws.Range("J66").Value = 34400
$ws.Range("L66").Value = 103200
$ws.Range("N66").Value = -110064

$ws.Range("H113").Value = 1238.125
$ws.Range("I113").Value = 985.46155
$ws.Range("K113").Value = 985.46155
$ws.Range("M113").Value = 1184.53845

$ws.Range("H126").Value = 2747.158
$ws.Range("I126").Value = 2014
$ws.Range("J126").Value = 3755.25
$ws.Range("K126").Value = 6042
$ws.Range("L126").Value = 11265.75
$ws.Range("M126").Value = -3572
$ws.Range("N126").Value = -16205.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 42084.2
$ws.Range("I34").Value = 42084.2
$ws.Range("K34").Value = 42084.2
$ws.Range("M34").Value = -41912.2

$ws.Range("H92").Value = 32000
$ws.Range("J92").Value = 32000
$ws.Range("L92").Value = 32000
$ws.Range("N92").Value = -36992

$ws.Range("H103").Value = 42200.668
$ws.Range("J103").Value = 42200.668
$ws.Range("L103").Value = 42200.668
$ws.Range("N103").Value = -44544.668

$ws.Range("H104").Value = 23092.5
$ws.Range("J104").Value = 23092.5
$ws.Range("L104").Value = 23092.5
$ws.Range("N104").Value = -30080.5

$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("N119").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 10026
$ws.Range("I32").Value = 10026
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 10026
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -9709
$ws.Range("N32").ClearContents()

$ws.Range("H64").Value = 23113.5
$ws.Range("J64").Value = 23113.5
$ws.Range("L64").Value = 23113.5
$ws.Range("N64").Value = -23609.5

$ws.Range("H67").Value = 23113.5
$ws.Range("J67").Value = 23113.5
$ws.Range("L67").Value = 23113.5
$ws.Range("N67").Value = -24829.5

$ws.Range("H126").Value = 6313.643
$ws.Range("I126").Value = 7240.9165
$ws.Range("J126").Value = 750
$ws.Range("K126").Value = 21722.7495
$ws.Range("L126").Value = 2250
$ws.Range("M126").Value = -19252.7495
$ws.Range("N126").Value = -7190

$ws.Range("H136").Value = 929.6177
$ws.Range("I136").Value = 905.56665
$ws.Range("K136").Value = 2716.69995
$ws.Range("M136").Value = -166.6999500000002
